# Added Week 15 simulations: a new player, "N.Bellore", joins the roster
# right after "A.Peterson". This pushes every player that used to sit in
# columns K:T (D.Metcalf .. C.Parkinson) one column to the right (L:U) on
# both the "Rushing" and "Receiving" sheets, so the header row (and the
# "n" placeholder row beneath it) are rewritten in roster order to land
# the new player in the correct spot and grow the table by one column.

$wb = $excel.ActiveWorkbook

$players = @(
    "R.Wilson", "G.Smith", "J.Luton", "C.Carson", "R.Penny", "A.Collins",
    "T.Homer", "D.Dallas", "A.Peterson", "N.Bellore", "D.Metcalf",
    "T.Lockett", "F.Swain", "D.Eskridge", "P.Hart", "J.Ursua", "G.Everett",
    "W.Dissly", "T.Mabry", "C.Parkinson"
)

foreach ($sheetName in @("Rushing", "Receiving")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Grow the table by one column (U). Give the new header cell the same
    # formatting (bold font + border + center/top alignment) as the
    # existing last header cell (T1) before T1's own value is rewritten.
    $ws.Range("U1").Value = "placeholder"
    $ws.Range("T1").Copy()
    $ws.Range("U1").PasteSpecial(-4122)  # xlPasteFormats

    $col = 2  # column B
    foreach ($player in $players) {
        $ws.Cells.Item(1, $col).Value = $player
        $ws.Cells.Item(2, $col).Value = "n"
        $col = $col + 1
    }

    $ws.Range("A2").Value = "Yards list"
}
